$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New verb group ("uyandirmak" - to wake up - and related forms) appended
# after the existing vocabulary table. Arabic forms go in column A, Turkish
# meanings in column B, rows 33-39.
$turkish = @(
    "Uyandırdı",
    "Uyandırıyor",
    "Uyandır ! // İkaz et !",
    "Uyandırmak",
    "Varlıklı ve müreffeh oldu",
    "Kesin olarak bildi",
    "Kuruttu"
)
$arabic = @(
    "أَيْقَظَ",
    "يُوقِظُ",
    "أَيْقِظْ",
    "إِيقَاظٌ",
    "أَيْسَرَ",
    "أَيْقَنَ",
    "أَيْبَسَ"
)

$startRow = 33
$lastExistingRow = 35

# Prepare the brand-new rows (36-39) so they inherit the same look as the
# already-existing placeholder rows 33-35 (row height 69.75, column A style
# matching the rest of the table).
for ($row = $lastExistingRow + 1; $row -lt ($startRow + $turkish.Length); $row++) {
    $ws.Range("A$startRow").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Rows("${row}:${row}").RowHeight = 69.75
}
$excel.CutCopyMode = $false

# Write column B (Turkish) first, then column A (Arabic) - matches the order
# the values were entered in the source document.
for ($i = 0; $i -lt $turkish.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $turkish[$i]
}
for ($i = 0; $i -lt $arabic.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $arabic[$i]
}

# Restore view state: scrolled down so the new rows are visible, cursor left
# on the next empty row ready for further entries.
$ws.Activate()
$ws.Range("B42").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
